# Apply edits to the bread_coop_2023-01-09 workbook:
# 1. Update the "timestamp" column (O) for all data rows (2-398) from
#    "2023-01-09 14:17:57" to "2023-01-09 15:13:39".
# 2. Update a few "ratingAmount" values (column D) that changed.
# 3. Update the productAriaLabel text for row 112 (column M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-09 15:13:39"

$lastRow = 398

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp  # column O = timestamp
}

# ratingAmount (column D) updates
$ws.Cells.Item(60, 4).Value = 5
$ws.Cells.Item(78, 4).Value = 83
$ws.Cells.Item(90, 4).Value = 6
$ws.Cells.Item(116, 4).Value = 2

# productAriaLabel (column M) update for row 112
$ws.Cells.Item(112, 13).Value = "Prix Garantie Zwieback - Online kein Bestand 1.80 Schweizer Franken"
